# advcalc.docx edit: "add info on button"
#  1. Version 1.3 -> Version 1.4
#  2. New paragraph after "It takes two numbers..." with text about
#     clicking a button to see the answers.
#  3. The "_GoBack" bookmark (which Word keeps at the location of the most
#     recent edit) moves from the end of the "Power" paragraph to the end
#     of the new paragraph.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1. Bump the version number: "3" -> "4"
# ---------------------------------------------------------------------
$d.Content.Find.Execute("Version 1.3", $true, $false, $false, $false,
                         $false, $true, 1, $false, "Version 1.4", 2) | Out-Null

# ---------------------------------------------------------------------
# 2. Insert the new paragraph right after the paragraph that talks
#    about selecting the operation from a drop-down.
# ---------------------------------------------------------------------
$dropdownPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*select the operation from a drop-down.*") {
        $dropdownPara = $p
        break
    }
}

$dropdownPara.Range.InsertParagraphAfter()

# Find the freshly-created (empty) paragraph and fill it in.
$newPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Start -eq $dropdownPara.Range.End) {
        $newPara = $p
        break
    }
}
$insertionPoint = $d.Range($newPara.Range.Start, $newPara.Range.Start)
$insertionPoint.InsertAfter("The user then clicks a button to see the answers.")

# ---------------------------------------------------------------------
# 3. Move the "_GoBack" bookmark from the "Power" paragraph to the end
#    of the paragraph we just added (right before its paragraph mark).
# ---------------------------------------------------------------------
$oldBookmark = $d.Bookmarks("_GoBack")
$oldBookmark.Delete()

# Re-resolve the new paragraph (its end position shifted once text was
# inserted into it above).
$targetPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "The user then clicks a button to see the answers.*") {
        $targetPara = $p
        break
    }
}
$bookmarkPos = $targetPara.Range.End - 1

# NOTE: adding a bookmark with Bookmarks.Add at a collapsed range that
# sits exactly one character before a paragraph mark is unreliable in
# this runtime, so a throw-away character is inserted right after the
# desired position first (making it a normal, interior position), the
# bookmark is added there, and then the throw-away character is removed
# again. The bookmark stays correctly anchored once the character that
# follows it is deleted.
$placeholder = $d.Range($bookmarkPos, $bookmarkPos)
$placeholder.InsertAfter("X")

$bookmarkRange = $d.Range($bookmarkPos, $bookmarkPos)
$d.Bookmarks.Add("_GoBack", $bookmarkRange)

$d.Range($bookmarkPos, $bookmarkPos + 1).Delete()
